$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-05-18 to 2021-05-19
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-56
$ws.Range("D2").Value = 0.02406750690098922
$ws.Range("E2").Value = 0.0008486562942009446
$ws.Range("D3").Value = 0.0179868748685138
$ws.Range("E3").Value = -0.004282655246252709
$ws.Range("D4").Value = 0.01791801330098137
$ws.Range("E4").Value = -0.002866076081292412
$ws.Range("D5").Value = 0.02024355013671647
$ws.Range("E5").Value = 0.0002978850163837254
$ws.Range("D6").Value = 0.01945864498136798
$ws.Range("E6").Value = -0.001679462571976931
$ws.Range("D7").Value = 0.02703108311951071
$ws.Range("E7").Value = -0.006944444444444531
$ws.Range("D8").Value = 0.01920498542186998
$ws.Range("E8").Value = -0.001661129568106268
$ws.Range("D9").Value = 0.01965472537705356
$ws.Range("E9").Value = -0.0011876484560569
$ws.Range("D10").Value = 0.01904080699532373
$ws.Range("E10").Value = -0.0061297045482408
$ws.Range("D11").Value = 0.01954968230793629
$ws.Range("E11").Value = 0.0002985074626864481
$ws.Range("D12").Value = 0.01932208899151554
$ws.Range("E12").Value = 0
$ws.Range("D13").Value = 0.01969596450789219
$ws.Range("E13").Value = 0.007723304231027628
$ws.Range("D14").Value = 0.01869844439968226
$ws.Range("E14").Value = -0.003786775415088761
$ws.Range("D15").Value = 0.01729787014478534
$ws.Range("E15").Value = -0.0009896091044037103
$ws.Range("D16").Value = 0.01779779734410271
$ws.Range("E16").Value = -0.004546746234725729
$ws.Range("D17").Value = 0.01612313848849054
$ws.Range("E17").Value = -0.0006635700066356387
$ws.Range("D18").Value = 0.01497019354282747
$ws.Range("E18").Value = -0.02630005977286332
$ws.Range("D19").Value = 0.01649973734369616
$ws.Range("E19").Value = -0.0112825833225263
$ws.Range("D20").Value = 0.01848368968059807
$ws.Range("E20").Value = 0.01810145232582627
$ws.Range("D21").Value = 0.01960453813291975
$ws.Range("E21").Value = -0.005715306304697299
$ws.Range("D22").Value = 0.0210175619367491
$ws.Range("E22").Value = -0.002471169686985242
$ws.Range("D23").Value = 0.01937169266304313
$ws.Range("E23").Value = -0.006175628859767968
$ws.Range("D24").Value = 0.02103779245376427
$ws.Range("E24").Value = -0.01266759130836803
$ws.Range("D25").Value = 0.02093702891709252
$ws.Range("E25").Value = -0.01382488479262667
$ws.Range("D26").Value = 0.01943646922233212
$ws.Range("E26").Value = -0.004984086951300148
$ws.Range("D27").Value = 0.02012430680084816
$ws.Range("E27").Value = -0.02934637616718538
$ws.Range("D28").Value = 0.02725634214550663
$ws.Range("E28").Value = 0.01554404145077704
$ws.Range("D29").Value = 0.01945631069094316
$ws.Range("E29").Value = -0.0350929814037193
$ws.Range("D30").Value = 0.01268706298314341
$ws.Range("E30").Value = -0.0139525612916086
$ws.Range("D31").Value = 0.009021059870950698
$ws.Range("E31").Value = 0.01584905660377367
$ws.Range("D32").Value = 0.01670671109469759
$ws.Range("E32").Value = -0.003108808290155474
$ws.Range("D33").Value = 0.02017011725043542
$ws.Range("E33").Value = -0.03873584114109918
$ws.Range("D34").Value = 0.01802597423312967
$ws.Range("E34").Value = 0.02730206006453217
$ws.Range("D35").Value = 0.01820941055568075
$ws.Range("E35").Value = -0.02845849802371547
$ws.Range("D36").Value = 0.01737918126125019
$ws.Range("E36").Value = -0.01110340041637747
$ws.Range("D37").Value = 0.01951408437895766
$ws.Range("E37").Value = 0.007625826131164137
$ws.Range("D38").Value = 0.01978213872940876
$ws.Range("E38").Value = -0.001622498647917636
$ws.Range("D39").Value = 0.02533794446470198
$ws.Range("E39").Value = -0.01240633828767956
$ws.Range("D40").Value = 0.01781452642548065
$ws.Range("E40").Value = -0.02738589211618259
$ws.Range("D41").Value = 0.02267529718678126
$ws.Range("E41").Value = -0.03795209663029309
$ws.Range("D42").Value = 0.01941118107606315
$ws.Range("E42").Value = -0.01370906321401388
$ws.Range("D43").Value = 0.02001245538465848
$ws.Range("E43").Value = -0.01311249137336079
$ws.Range("D44").Value = 0.01847337989788841
$ws.Range("E44").Value = -0.02688302252361341
$ws.Range("D45").Value = 0.0202791480656951
$ws.Range("E45").Value = -0.005515587529976052
$ws.Range("D46").Value = 0.0196858492493846
$ws.Range("E46").Value = -0.02239130434782621
$ws.Range("D47").Value = 0.01783514599089996
$ws.Range("E47").Value = 0.001603298213467719
$ws.Range("D48").Value = 0.01627000426105265
$ws.Range("E48").Value = -0.02008608321377325
$ws.Range("D49").Value = 0.01727433271633498
$ws.Range("E49").Value = 0.01787101787101797
$ws.Range("D50").Value = 0.01749317244366263
$ws.Range("E50").Value = -0.03883106485188148
$ws.Range("D51").Value = 0.01617429835363469
$ws.Range("E51").Value = 0.02416173570019708
$ws.Range("D52").Value = 0.01803239353179795
$ws.Range("E52").Value = 0.001725997842502824
$ws.Range("D53").Value = 0.01565491873411042
$ws.Range("E53").Value = -0.009716941275876745
$ws.Range("D54").Value = 0.007513497304915811
$ws.Range("E54").Value = -0.0003883495145632354
$ws.Range("D55").Value = 0.007279873738230921
$ws.Range("E55").Value = -0.01239846088071817
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = -0.00724962522480932

# Restore sheet protection
$ws.Protect()

